$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 18, shifting the existing
# rows 18-22 down to 20-24 (all their data/styles move with them).
$ws.Rows("18:19").Insert()

# Row 18: new weekly record (Primera, bandeja 10 kilos)
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 45126
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100107
$ws.Range("H18").Value = "Otros"
$ws.Range("I18").Value = 100107002
$ws.Range("J18").Value = "Chirimoya"
$ws.Range("K18").Value = "Cultivar IV Región"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 14000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 14375
$ws.Range("Q18").Value = "$/bandeja 10 kilos"
$ws.Range("R18").Value = "Región de Coquimbo"
$ws.Range("S18").Value = 1438
$ws.Range("T18").Value = 10

# Row 19: new weekly record (Segunda, bandeja 10 kilos)
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 45126
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100107
$ws.Range("H19").Value = "Otros"
$ws.Range("I19").Value = 100107002
$ws.Range("J19").Value = "Chirimoya"
$ws.Range("K19").Value = "Cultivar IV Región"
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 180
$ws.Range("N19").Value = 13000
$ws.Range("O19").Value = 13000
$ws.Range("P19").Value = 13000
$ws.Range("Q19").Value = "$/bandeja 10 kilos"
$ws.Range("R19").Value = "Región de Coquimbo"
$ws.Range("S19").Value = 1300
$ws.Range("T19").Value = 10
